$d = $word.ActiveDocument

# Each bullet/paragraph is replaced wholesale via Range.InsertXML so the
# resulting markup exactly matches a clean English paragraph (single run,
# no leftover run-level rPr) instead of leaving old Spanish runs behind.

$xml1 = '<w:p><w:r><w:t>ContosoLearn Market Research</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML($xml1)

$xml2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology.</w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($xml2)

$xml3 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations.</w:t></w:r></w:p>'
$d.Paragraphs(3).Range.InsertXML($xml3)

$xml4 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning.</w:t></w:r></w:p>'
$d.Paragraphs(4).Range.InsertXML($xml4)

$xml5 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration.</w:t></w:r></w:p>'
$d.Paragraphs(5).Range.InsertXML($xml5)

$xml6 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs.</w:t></w:r></w:p>'
$d.Paragraphs(6).Range.InsertXML($xml6)

$xml7 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario.</w:t></w:r></w:p>'
$d.Paragraphs(7).Range.InsertXML($xml7)

$xml8 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>Munson''sLearn: Munson''sLearn is designed to enable businesses to train their employees, partners, and customers.</w:t></w:r></w:p>'
$d.Paragraphs(8).Range.InsertXML($xml8)

$xml9 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project.</w:t></w:r></w:p>'
$d.Paragraphs(9).Range.InsertXML($xml9)

$xml10 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a best</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-in-class training experience.</w:t></w:r></w:p>'
$d.Paragraphs(10).Range.InsertXML($xml10)

$xml11 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises.</w:t></w:r></w:p>'
$d.Paragraphs(11).Range.InsertXML($xml11)

$xml12 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website.</w:t></w:r></w:p>'
$d.Paragraphs(12).Range.InsertXML($xml12)

$xml13 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects.</w:t></w:r></w:p>'
$d.Paragraphs(13).Range.InsertXML($xml13)

$xml14 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010.</w:t></w:r></w:p>'
$d.Paragraphs(14).Range.InsertXML($xml14)

$xml15 = '<w:p><w:r><w:t xml:space="preserve">These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. </w:t></w:r></w:p>'
$d.Paragraphs(15).Range.InsertXML($xml15)
